# Add "supplier_id" and "stok" columns (F and G) to the barang template,
# matching the "Add stok and penjualan management features" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same (bold) as the existing headers in row 1.
$ws.Range("F1").Value = "supplier_id"
$ws.Range("G1").Value = "stok"
$ws.Range("F1:G1").Font.Bold = $true

# supplier_id values for each existing data row.
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 2

# stok values for each existing data row.
$ws.Range("G2").Value = 10
$ws.Range("G3").Value = 20
$ws.Range("G4").Value = 30
$ws.Range("G5").Value = 40
$ws.Range("G6").Value = 50

# Update the active selection to reflect where the author left the cursor.
$ws.Range("G14").Select()
